# "Se libera version 1.1.0 de PC y 1.03 de PLC"
#
# Populates the "Hoja2" sheet with the new PLC register table
# (H2-01, A1-03, B1-01, B1-02, E1-09) and updates the saved
# selection / active-sheet state so that Hoja2 becomes the
# active tab with C7 selected, while Hoja1 keeps B8 selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- New data on Hoja2 ------------------------------------------------
$ws2.Range("B2").Value = "H2-01"
$ws2.Range("C2").Value = 2220

$ws2.Range("B3").Value = "A1-03"
$ws2.Range("C3").Value = 1

$ws2.Range("B4").Value = "B1-01"
$ws2.Range("C4").Value = 2

$ws2.Range("B5").Value = "B1-02"
$ws2.Range("C5").Value = 2

$ws2.Range("B6").Value = "E1-09"
$ws2.Range("C6").Value = 0.01

# --- Selection / active sheet state -----------------------------------
# Leave Hoja1 selected on B8 (no longer the active tab).
$null = $ws1.Activate()
$null = $ws1.Range("B8").Select()

# Make Hoja2 the active tab, with C7 selected.
$null = $ws2.Activate()
$null = $ws2.Range("C7").Select()
